$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime (H2) and
# Correspond Handback DateTime (K2) for the 051d8a20... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-29 08:51:59"
$wsZhCn.Range("K2").Value = "2016-08-29 08:52:29"

# de-de sheet: update Correspond Handback DateTime (K2) for the
# 051d8a20... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-29 08:52:36"
